$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the "Files" Bento query (cell B4): drop the `File Type` and
# `Breed` coalesce lines from the RETURN clause.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Pulmonary Adenocarcinoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newQuery

# Move the active selection from E4 to B4, matching the saved view state.
$ws.Range("B4").Select()
